# Updates cryptos price/volume columns (D, E) for Sheet1, row by row,
# matching the scraped values from the latest GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.764.33'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '2.088.37'
$ws.Range('E3').Value = '  +0.92%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.98'
$ws.Range('E5').Value = '  +0.47%  '
$ws.Range('E6').Value = '  +2.45%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '58.20'
$ws.Range('E8').Value = '  +0.70%  '
$ws.Range('E9').Value = '  +1.04%  '
$ws.Range('E10').Value = '  +0.17%  '
$ws.Range('E11').Value = '  +3.10%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.25'
$ws.Range('E12').Value = '  +3.51%  '
$ws.Range('D13').Value = '2.395.69'
$ws.Range('E13').Value = '  +0.83%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.11'
$ws.Range('E14').Value = '  +1.69%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.777'
$ws.Range('E15').Value = '  +1.13%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.36'
$ws.Range('E16').Value = '  +1.86%  '
$ws.Range('D17').Value = '2.080.92'
$ws.Range('E17').Value = '  +0.18%  '
$ws.Range('D18').Value = '37.780.07'
$ws.Range('E18').Value = '  +0.57%  '
$ws.Range('E19').Value = '  -0.45%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.02'
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('E21').Value = '  +0.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '229.60'
$ws.Range('E22').Value = '  +1.09%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('E24').Value = '  -0.95%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.40'
$ws.Range('E25').Value = '  +1.13%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.71'
$ws.Range('E26').Value = '  +8.61%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '170.81'
$ws.Range('E27').Value = '  +0.93%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.132'
$ws.Range('E28').Value = '  -4.19%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.56'
$ws.Range('E29').Value = '  +1.11%  '
$ws.Range('E30').Value = '  +0.91%  '
$ws.Range('E31').Value = '  +2.28%  '
$ws.Range('E32').Value = '  +0.85%  '
$ws.Range('E33').Value = '  +1.38%  '
$ws.Range('E34').Value = '  -0.40%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.49'
$ws.Range('E35').Value = '  +1.93%  '
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('E37').Value = '  -1.07%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.39'
$ws.Range('E39').Value = '  +1.62%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0235'
$ws.Range('E40').Value = '  +9.94%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '101.22'
$ws.Range('E41').Value = '  +3.74%  '
$ws.Range('E42').Value = '  -0.95%  '
$ws.Range('E43').Value = '  +1.07%  '
$ws.Range('E44').Value = '  +3.46%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '16.75'
$ws.Range('E45').Value = '  +2.60%  '
$ws.Range('D46').Value = '1.452.31'
$ws.Range('E46').Value = '  -0.01%  '
$ws.Range('E47').Value = '  +0.17%  '
$ws.Range('E48').Value = '  -3.36%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.21'
$ws.Range('E49').Value = '  -1.62%  '
$ws.Range('E50').Value = '  -1.43%  '
$ws.Range('D51').Value = '2.279.13'
$ws.Range('E51').Value = '  +0.81%  '
